$d = $word.ActiveDocument

# The title paragraph currently reads "Sponsor’sship Levels: Bionic Badgers Of
# Wyoming FTC Team" (a typo: "Sponsor's" + "ship" glued together). Fix the
# typo so it reads "Sponsorship Levels: ..." by removing the stray
# apostrophe-s, without disturbing the run/bookmark layout of the rest of
# the paragraph.
$titlePara = $d.Paragraphs(1).Range

# Locate "Sponsor’s" within the title paragraph.
$aposRange = $titlePara.Duplicate
$aposRange.Find.Execute("Sponsor’s", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$wordStart = $aposRange.Start
$aposStart = $wordStart + 7
$aposEnd = $aposRange.End

# Locate the boundary right after "ship " (between "ship " and "Levels").
$shipRange = $titlePara.Duplicate
$shipRange.Find.Execute("ship ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$shipEnd = $shipRange.End

# Drop a temporary bookmark at the "ship "/"Levels" boundary so that when we
# delete the "’s" text below, Word's run-merge doesn't fuse those two runs
# together (a bookmark anchor forces a run split, just like the real
# "_GoBack" bookmark does below).
$guard = $d.Range($shipEnd, $shipEnd)
$d.Bookmarks.Add("ZZGuardBoundary", $guard)

# Move the "_GoBack" bookmark to sit right after "Sponsor" (i.e. where the
# user's last edit happened) — this both splits "Sponsor" from "’s" into
# separate runs and relocates the bookmark from its old spot further down
# in the document (Word bookmarks are unique by name, so re-adding moves
# it).
$editPoint = $d.Range($aposStart, $aposStart)
$d.Bookmarks.Add("_GoBack", $editPoint)

# Remove the stray "’s".
$aposTextRange = $d.Range($aposStart, $aposEnd)
$aposTextRange.Text = ""

# Drop the temporary guard bookmark now that the edit is done.
$d.Bookmarks("ZZGuardBoundary").Delete()
